$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 7-39 with the reshuffled workflow transition data
$ws.Cells.Item(7, 1).Value = "CompletamentoManualeIncarico"
$ws.Cells.Item(7, 2).Value = "VenditaRiparazione"
$ws.Cells.Item(7, 3).Value = "VenditaRiparazione"
$ws.Cells.Item(8, 1).Value = "CompletamentoManualeIncarico"
$ws.Cells.Item(8, 2).Value = "RiparazioneNonPrevista"
$ws.Cells.Item(8, 3).Value = "EsecuzioneTriage"
$ws.Cells.Item(9, 1).Value = "VenditaRiparazione"
$ws.Cells.Item(9, 2).Value = "RiparazioneVenduta"
$ws.Cells.Item(9, 3).Value = "AvvioRiparazione"
$ws.Cells.Item(10, 1).Value = "VenditaRiparazione"
$ws.Cells.Item(10, 2).Value = "CampiObbligMancanti"
$ws.Cells.Item(10, 3).Value = "VenditaRiparazione"
$ws.Cells.Item(11, 1).Value = "VenditaRiparazione"
$ws.Cells.Item(11, 2).Value = "RiparazioneNonVenduta"
$ws.Cells.Item(11, 3).Value = "EsecuzioneTriage"
$ws.Cells.Item(12, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(12, 2).Value = "SopralluogoPCE"
$ws.Cells.Item(12, 3).Value = "AvvioSopralluogoPCE"
$ws.Cells.Item(13, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(13, 2).Value = "Visio"
$ws.Cells.Item(13, 3).Value = "AvvioVisio"
$ws.Cells.Item(14, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(14, 2).Value = "Desk"
$ws.Cells.Item(14, 3).Value = "AvvioDesk"
$ws.Cells.Item(15, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(15, 2).Value = "Sopralluogo"
$ws.Cells.Item(15, 3).Value = "AvvioSopralluogo"
$ws.Cells.Item(16, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(16, 2).Value = "VenditaRiparazionePA1"
$ws.Cells.Item(16, 3).Value = "VenditaRiparazionePA1"
$ws.Cells.Item(17, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(17, 2).Value = "InvioLinkSelfCare"
$ws.Cells.Item(17, 3).Value = "InvioLinkSelfCare"
$ws.Cells.Item(18, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(18, 2).Value = "NessunServizio"
$ws.Cells.Item(18, 3).Value = "SceltaManualeServizio"
$ws.Cells.Item(19, 1).Value = "EsecuzioneTriage"
$ws.Cells.Item(19, 2).Value = "CampiObbligMancanti"
$ws.Cells.Item(19, 3).Value = "SceltaManualeServizio"
$ws.Cells.Item(20, 1).Value = "VenditaRiparazionePA1"
$ws.Cells.Item(20, 2).Value = "InteressatoRiparazione"
$ws.Cells.Item(20, 3).Value = "VenditaRiparazione"
$ws.Cells.Item(21, 1).Value = "VenditaRiparazionePA1"
$ws.Cells.Item(21, 2).Value = "NonInteressatoRiparazione"
$ws.Cells.Item(21, 3).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(22, 1).Value = "VenditaRiparazionePA1"
$ws.Cells.Item(22, 2).Value = "TimerScaduto"
$ws.Cells.Item(22, 3).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(23, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(23, 2).Value = "SopralluogoPCE"
$ws.Cells.Item(23, 3).Value = "AvvioSopralluogoPCE"
$ws.Cells.Item(24, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(24, 2).Value = "Visio"
$ws.Cells.Item(24, 3).Value = "AvvioVisio"
$ws.Cells.Item(25, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(25, 2).Value = "Desk"
$ws.Cells.Item(25, 3).Value = "AvvioDesk"
$ws.Cells.Item(26, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(26, 2).Value = "Sopralluogo"
$ws.Cells.Item(26, 3).Value = "AvvioSopralluogo"
$ws.Cells.Item(27, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(27, 2).Value = "SelfCare"
$ws.Cells.Item(27, 3).Value = "AvvioAccertSelfCare"
$ws.Cells.Item(28, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(28, 2).Value = "NessunServizio"
$ws.Cells.Item(28, 3).Value = "SceltaManualeServizio"
$ws.Cells.Item(29, 1).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(29, 2).Value = "CampiObbligMancanti"
$ws.Cells.Item(29, 3).Value = "SceltaManualeServizio"
$ws.Cells.Item(30, 1).Value = "InvioLinkSelfCare"
$ws.Cells.Item(30, 2).Value = "SelfCareEseguito"
$ws.Cells.Item(30, 3).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(31, 1).Value = "InvioLinkSelfCare"
$ws.Cells.Item(31, 2).Value = "SelfCareNonEseguito"
$ws.Cells.Item(31, 3).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(32, 1).Value = "InvioLinkSelfCare"
$ws.Cells.Item(32, 2).Value = "TimerScaduto"
$ws.Cells.Item(32, 3).Value = "SceltaAutomaticaServizio"
$ws.Cells.Item(33, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(33, 2).Value = "AvvioRiparazione"
$ws.Cells.Item(33, 3).Value = "AvvioRiparazione"
$ws.Cells.Item(34, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(34, 2).Value = "SopralluogoPCE"
$ws.Cells.Item(34, 3).Value = "AvvioSopralluogoPCE"
$ws.Cells.Item(35, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(35, 2).Value = "Sopralluogo"
$ws.Cells.Item(35, 3).Value = "AvvioSopralluogo"
$ws.Cells.Item(36, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(36, 2).Value = "Visio"
$ws.Cells.Item(36, 3).Value = "AvvioManualeVisio"
$ws.Cells.Item(37, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(37, 2).Value = "Desk"
$ws.Cells.Item(37, 3).Value = "AvvioManualeDesk"
$ws.Cells.Item(38, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(38, 2).Value = "InvioLinkSelfCare"
$ws.Cells.Item(38, 3).Value = "InvioLinkSelfCare"
$ws.Cells.Item(39, 1).Value = "SceltaManualeServizio"
$ws.Cells.Item(39, 2).Value = "CampiObbligMancanti"
$ws.Cells.Item(39, 3).Value = "SceltaManualeServizio"

# Remove the now-obsolete last row (table shrank from 40 to 39 rows)
$ws.Rows.Item(40).Delete()
